# "10 years Finalization data"
# Build a second worksheet ("Sheet1") that holds a clean copy of the daily
# data table (header row + 30 daily rows) that already lives on the
# "Data Harian - Table" sheet, and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# Make sure gridlines are shown on the source sheet (matches the workbook's
# on-disk state) and select the table that is about to be copied - this
# mirrors the manual "select A9:K39, copy" workflow the author used.
$wb.Windows.Item(1).DisplayGridlines = $true
[void]$dataSheet.Range("A9:K39").Select()
[void]$dataSheet.Range("A9:K39").Copy()

# Insert the new sheet right after "Data Harian - Table"; Excel names it
# "Sheet1" by default and activates it automatically.
$summarySheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)

# Paste values first, then formats, so the new sheet reuses the existing
# style records (header style + bordered data style) instead of minting new
# ones.
[void]$summarySheet.Range("A1").PasteSpecial(-4163)
[void]$dataSheet.Range("A9:K39").Copy()
[void]$summarySheet.Range("A1").PasteSpecial(-4122)

# Leave the new sheet selected over the whole pasted table, as the active tab.
[void]$summarySheet.Range("A1:K31").Select()
$summarySheet.Activate()
